# Generate Report for Handoff
# Adds a new handoff entry (a05cd29b-b1ea-4b69-b814-7251f5d09085) as a new
# row to the Overview, zh-cn and de-de worksheets, mirroring the existing
# 832cefc6-67c5-4d58-a591-a6122c5ce5fe row already present in each sheet.

$wb = $excel.ActiveWorkbook

# ---- shared literals -------------------------------------------------
$mdName        = "a05cd29b-b1ea-4b69-b814-7251f5d09085.md"
$zhXlfName     = "a05cd29b-b1ea-4b69-b814-7251f5d09085.f20852dc2714d4e947a96593badcffc723818fd0.zh-cn.xlf"
$deXlfName     = "a05cd29b-b1ea-4b69-b814-7251f5d09085.f20852dc2714d4e947a96593badcffc723818fd0.de-de.xlf"

$mdUrl    = "https://github.com/OpenLocalizationTest/oltest/blob/88328c67bb05bbc02c681f9bb925a4de657c7e76/e2e/$mdName"
$zhXlfUrl = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/a8b15527fdb0a1bb775c74021a466b23031c80e6/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/$zhXlfName"
$deXlfUrl = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/606d935de2bc320921a4977e312559fa388e550f/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/$deXlfName"

$readyStatus   = "Ready for handoff"
$handoffDate   = "2016-29-20 08:29:43"
$zhHandoffDT   = "2016-03-20 08:29:40"
$deHandoffDT   = "2016-03-20 08:29:43"
$epoch         = "0001-01-01 00:00:00"
$includeText   = "Include"
$dateFormat    = "yyyy-mm-dd HH:mm:ss"

# =======================================================================
# Sheet "Overview": File Name | zh-cn | de-de | Latest Handoff Date
# =======================================================================
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Hyperlinks.Add($wsOverview.Range("A3"), $mdUrl, "", "", $mdName)
$wsOverview.Range("B3").Value = $readyStatus
$wsOverview.Range("C3").Value = $readyStatus
$wsOverview.Range("D3").Value = $handoffDate

# =======================================================================
# Sheet "zh-cn": Source File Name | File Extension | Status |
#   Latest Handoff File | Latest Handoff Datetime | Latest Target File |
#   Latest Handback File | Latest Handback DateTime | Handoff Reason |
#   Dependency From | Error Detail
# =======================================================================
$wsZhCn = $wb.Worksheets.Item("zh-cn")

$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A3"), $mdUrl, "", "", $mdName)
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("B3"), $mdUrl, "", "", ".md")
$wsZhCn.Range("C3").Value = $readyStatus
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("D3"), $zhXlfUrl, "", "", $zhXlfName)
$wsZhCn.Range("E3").Value = $zhHandoffDT
$wsZhCn.Range("E3").NumberFormat = $dateFormat
$wsZhCn.Range("H3").Value = $epoch
$wsZhCn.Range("I3").Value = $includeText

# =======================================================================
# Sheet "de-de": same layout as zh-cn
# =======================================================================
$wsDeDe = $wb.Worksheets.Item("de-de")

$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A3"), $mdUrl, "", "", $mdName)
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("B3"), $mdUrl, "", "", ".md")
$wsDeDe.Range("C3").Value = $readyStatus
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("D3"), $deXlfUrl, "", "", $deXlfName)
$wsDeDe.Range("E3").Value = $deHandoffDT
$wsDeDe.Range("E3").NumberFormat = $dateFormat
$wsDeDe.Range("H3").Value = $epoch
$wsDeDe.Range("I3").Value = $includeText
